# Add "Test" to cell B5 and leave it selected, matching the upload edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Test"
$ws.Range("B5").Select()
